$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "2025-11-30T11:25:30.596Z"
$ws.Range("B35").Value = "rousol@testhost.com"
$ws.Range("C35").Value = "https://dga.gov.sa"
$ws.Range("D35").Value = "Yes"
$ws.Range("E35").Value = "٣٠‏/١١‏/٢٠٢٥"
$ws.Range("F35").Value = "٢:٢٥:٣٠ م"

$ws.Range("A36").Value = "2025-11-30T11:33:22.258Z"
$ws.Range("B36").Value = "sara@test.com"
$ws.Range("C36").Value = "https://www.arabou.edu.sa/"
$ws.Range("D36").Value = "Yes"
$ws.Range("E36").Value = "٣٠‏/١١‏/٢٠٢٥"
$ws.Range("F36").Value = "٢:٣٣:٢٢ م"
